$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Reorder the tied_teams (column O) string lists per the updated UEFA
# suspense rule: in a head-to-head situation the goal-diff gap can be
# -2 now (not just -1), which changes the ordering of tied teams shown.

for ($r = 40; $r -le 52; $r++) {
    $ws.Cells.Item($r, 15).Value = "['Ireland', 'Costa Rica']"
}

for ($r = 53; $r -le 59; $r++) {
    $ws.Cells.Item($r, 15).Value = "['Argentina', 'Ireland', 'Colombia', 'Costa Rica']"
}

for ($r = 60; $r -le 62; $r++) {
    $ws.Cells.Item($r, 15).Value = "['Argentina', 'Colombia']"
}

for ($r = 63; $r -le 73; $r++) {
    $ws.Cells.Item($r, 15).Value = "['Argentina', 'Colombia', 'Scotland', 'Austria']"
}

$ws.Cells.Item(78, 15).Value = "['South Korea', 'Netherlands']"

for ($r = 101; $r -le 102; $r++) {
    $ws.Cells.Item($r, 15).Value = "['Italy', 'Netherlands']"
}
